# Updated cryptos list on Wed Jan 17 06:52:20 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) for each coin row,
# plus the two coin-pair row swaps (ShibaInu/InternetComputer and
# NEARProtocol/Maker). Numeric-looking prices are entered with a leading
# apostrophe so Excel keeps them as literal text (matching the source
# data, which stores prices like "42.848.99" / "310.98" as strings, not
# numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.848.99'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.565.12'
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''310.98'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").Value = '''98.75'
$ws.Range("E6").Value = '  +3.21%  '
$ws.Range("E7").Value = '  -0.40%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +0.12%  '
$ws.Range("D10").Value = '''35.90'
$ws.Range("E10").Value = '  +0.41%  '
$ws.Range("D11").Value = '''0.0809'
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '''7.47'
$ws.Range("E12").Value = '  -0.71%  '
$ws.Range("D13").Value = '2.958.77'
$ws.Range("E13").Value = '  +1.34%  '
$ws.Range("E14").Value = '  -1.69%  '
$ws.Range("D15").Value = '''15.84'
$ws.Range("E15").Value = '  +5.23%  '
$ws.Range("D16").Value = '2.565.99'
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").Value = '''0.842'
$ws.Range("E17").Value = '  -0.68%  '
$ws.Range("D18").Value = '42.864.32'
$ws.Range("E18").Value = '  +0.14%  '
$ws.Range("D19").Value = '''6.72'
$ws.Range("E19").Value = '  -1.15%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.0₃0962'
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").Value = '''12.42'
$ws.Range("E21").Value = '  -1.89%  '
$ws.Range("D22").Value = '''69.49'
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("D23").Value = '''248.30'
$ws.Range("E23").Value = '  -1.23%  '
$ws.Range("D24").Value = '''2.92'
$ws.Range("E24").Value = '  -0.88%  '
$ws.Range("E25").Value = '  +0.63%  '
$ws.Range("D26").Value = '''27.11'
$ws.Range("E26").Value = '  +2.34%  '
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("D29").Value = '''40.03'
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("D30").Value = '''10.20'
$ws.Range("E30").Value = '  -2.02%  '
$ws.Range("D31").Value = '''159.09'
$ws.Range("D32").Value = '''5.78'
$ws.Range("E32").Value = '  -2.15%  '
$ws.Range("D33").Value = '''0.0799'
$ws.Range("E33").Value = '  +2.61%  '
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("E35").Value = '  -1.75%  '
$ws.Range("D36").Value = '''3.29'
$ws.Range("E36").Value = '  -0.56%  '
$ws.Range("D37").Value = '''18.71'
$ws.Range("E37").Value = '  -1.13%  '
$ws.Range("E38").Value = '  +12.61%  '
$ws.Range("D39").Value = '''0.112'
$ws.Range("E39").Value = '  +0.85%  '
$ws.Range("E40").Value = '  -0.15%  '
$ws.Range("D41").Value = '''22.88'
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("D42").Value = '''4.12'
$ws.Range("E42").Value = '  +7.90%  '
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").Value = '''0.0302'
$ws.Range("E44").Value = '  -0.96%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = '''3.22'
$ws.Range("E45").Value = '  -1.01%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.992.44'
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '''9.05'
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("D48").Value = '2.810.76'
$ws.Range("E48").Value = '  +1.29%  '
$ws.Range("D49").Value = '''0.195'
$ws.Range("E49").Value = '  +2.81%  '
$ws.Range("D50").Value = '''81.39'
$ws.Range("E50").Value = '  -3.54%  '
$ws.Range("D51").Value = '''74.05'
